$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "0.0003800392150878906   1478.7844496356363"
$ws.Range("C3").Value = "0.0011472702026367188   1586.6035551964073"
$ws.Range("C4").Value = "0.0023202896118164062   2700.169578241978"
$ws.Range("C5").Value = "0.003298044204711914   2304.1943061377356"
$ws.Range("C6").Value = "0.026426076889038086   2118.735772259764"
$ws.Range("C7").Value = "0.12920498847961426   3115.6124635030487"
$ws.Range("C8").Value = "1.1275792121887207   3125.1921031396832"
$ws.Range("C9").Value = "12.294594049453735   3058.6771489783323"
$ws.Range("C10").Value = "149.77348399162292   2849.2265318034906"
